$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Resize/move the subtitle placeholder (shape 2) to make room for the new textbox ---
$subtitle = $s.Shapes.Item(2)
$subtitle.Left = 406.2856602913386
$subtitle.Top = 326.8173985748032
$subtitle.Width = 262.2087251574803
$subtitle.Height = 64.94795675590552

# --- 2. Add the new "Course Material" textbox ---
$tb = $s.Shapes.AddTextbox(1, 122.67976077952756, 350.5669411338583, 237.32024422047246, 41.1984241968504)
$tb.Name = "TextBox 2"
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

# First paragraph
$tb.TextFrame.TextRange.Text = "Course Material:"
$tb.TextFrame.TextRange.Font.Name = "Avenir Book"
$tb.TextFrame.TextRange.Font.Color.RGB = 5984069

# Second paragraph, appended after the first so both keep their own run identity
$para2 = $tb.TextFrame.TextRange.InsertAfter([char]13 + "https://github.com/CWML/Python1")
$para2.Font.Name = "Avenir Book"
$para2.Font.Color.RGB = 5984069

# Split the second paragraph into three runs: "https://" / "github.com" / "/CWML/Python1"
$secondParagraph = $tb.TextFrame.TextRange.Paragraphs(2)
$start = $secondParagraph.Start
$run1 = $tb.TextFrame.TextRange.Characters($start, 8)
$run2 = $tb.TextFrame.TextRange.Characters($start + 8, 10)
$run3 = $tb.TextFrame.TextRange.Characters($start + 18, 13)
$run1.Font.Name = "Avenir Book"
$run1.Font.Color.RGB = 5984069
$run2.Font.Name = "Avenir Book"
$run2.Font.Color.RGB = 5984069
$run3.Font.Name = "Avenir Book"
$run3.Font.Color.RGB = 5984069

# Restore the exact autofit size computed from the two lines of text
$tb.Height = 41.1984241968504
